# Generate Report for handback
#
# Adds three new source-file entries to the report:
#   fffffff61273d2-8bf6-43d7-b3e2-128ffd4c5234.md  (row 4, reuses dedb3467's handoff/handback info)
#   2f896ebd-0648-4060-aacb-62692ef7c544.md         (row 5)
#   6329960e-3640-4a37-8321-56d9dcf5ee52.md         (row 6)
# across the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText  = "Handed back: in sync with en-US"
$includeText = "Include"

# Commit SHA used for the existing "e2e/<file>.md" source links.
$srcCommit = "593c9aff7438ff881575552994205808e22e177c"

$mdA = "fffffff61273d2-8bf6-43d7-b3e2-128ffd4c5234.md"
$mdB = "2f896ebd-0648-4060-aacb-62692ef7c544.md"
$mdC = "6329960e-3640-4a37-8321-56d9dcf5ee52.md"

function Style-LikeHyperlink($cell) {
    # Re-create the look of the workbook's custom "HyperLink" cell style
    # (underlined, cornflowerblue) for cells that carry a hyperlink.
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}

function Add-Link($ws, $cell, $url, $display) {
    $ws.Hyperlinks.Add($cell, $url, "", "", $display) | Out-Null
    Style-LikeHyperlink $cell
}

# ---------------------------------------------------------------------------
# 1. Overview sheet - rows 4, 5, 6
# ---------------------------------------------------------------------------
$overviewRows = @(
    @{ Row = 4; Md = $mdA },
    @{ Row = 5; Md = $mdB },
    @{ Row = 6; Md = $mdC }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $wsOverview.Cells.Item($row, 1).Value2 = $r.Md
    $wsOverview.Cells.Item($row, 2).Value2 = $statusText
    $wsOverview.Cells.Item($row, 3).Value2 = $statusText

    $url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$($r.Md)"
    Add-Link $wsOverview $wsOverview.Cells.Item($row, 1) $url $r.Md
}

# ---------------------------------------------------------------------------
# 2. zh-cn / de-de sheet - row 4 (re-uses the dedb3467 handoff/handback data)
# ---------------------------------------------------------------------------
$dedbMd      = "dedb3467-9493-4ab0-8961-60a966699708.md"
$dedbZhXlf   = "dedb3467-9493-4ab0-8961-60a966699708.ea5a5ce460d2cdef84e4e3672d5db72352d42d69.zh-cn.xlf"
$dedbDeXlf   = "dedb3467-9493-4ab0-8961-60a966699708.ea5a5ce460d2cdef84e4e3672d5db72352d42d69.de-de.xlf"

function Fill-Row($ws, $row, $mdName, $xlfName, $hoffDt, $srcMdName, $backDt) {
    $ws.Cells.Item($row, 1).Value2 = $mdName
    $ws.Cells.Item($row, 2).Value2 = $statusText
    $ws.Cells.Item($row, 3).Value2 = $xlfName
    $ws.Cells.Item($row, 4).Value2 = $hoffDt
    $ws.Cells.Item($row, 5).Value2 = $srcMdName
    $ws.Cells.Item($row, 6).Value2 = $xlfName
    $ws.Cells.Item($row, 7).Value2 = $backDt
    $ws.Cells.Item($row, 8).Value2 = $includeText
}

function Link-Row($ws, $row, $mdUrl, $mdName, $xlfUrl, $xlfName, $srcUrl, $srcMdName) {
    Add-Link $ws $ws.Cells.Item($row, 1) $mdUrl  $mdName
    Add-Link $ws $ws.Cells.Item($row, 3) $xlfUrl $xlfName
    Add-Link $ws $ws.Cells.Item($row, 5) $srcUrl $srcMdName
    Add-Link $ws $ws.Cells.Item($row, 6) $xlfUrl $xlfName
}

$zhRow4Url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdA"
$zhXlf4Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/afcfafd7368a4fe5350019f2e3676d53b37ba336/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$dedbZhXlf"
$zhSrc4Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4f31164c5f3306f2d6f1b09a5c6446aa383d4c79/e2e/$dedbMd"

Fill-Row $wsZhCn 4 $mdA $dedbZhXlf "2016-01-25 09:14:25" $dedbMd "2016-01-25 09:15:22"
Link-Row $wsZhCn 4 $zhRow4Url $mdA $zhXlf4Url $dedbZhXlf $zhSrc4Url $dedbMd

$deRow4Url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdA"
$deXlf4Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74283cc9a8ec6c1c113946216d22838362a1e962/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$dedbDeXlf"
$deSrc4Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e95e331e240825b31d89d4f4acc5965a80537a70/e2e/$dedbMd"

Fill-Row $wsDeDe 4 $mdA $dedbDeXlf "2016-01-25 09:14:37" $dedbMd "2016-01-25 09:15:40"
Link-Row $wsDeDe 4 $deRow4Url $mdA $deXlf4Url $dedbDeXlf $deSrc4Url $dedbMd

# ---------------------------------------------------------------------------
# 3. zh-cn sheet - rows 5 and 6 (new handoff/handback info)
# ---------------------------------------------------------------------------
$zhXlf5 = "2f896ebd-0648-4060-aacb-62692ef7c544.9e4bc652c0f3b7031a37aa768eb56a3aab3d784e.zh-cn.xlf"
$zhXlf6 = "6329960e-3640-4a37-8321-56d9dcf5ee52.2daef6e2b6c383f8d79ffb9a56a7f70dbdb68b3b.zh-cn.xlf"

Fill-Row $wsZhCn 5 $mdB $zhXlf5 "2016-01-25 09:18:37" $mdB "2016-01-25 09:19:23"
Fill-Row $wsZhCn 6 $mdC $zhXlf6 "2016-01-25 09:18:37" $mdC "2016-01-25 09:19:23"

$zhXlf5Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/afcfafd7368a4fe5350019f2e3676d53b37ba336/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$zhXlf5"
$zhSrc5Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4f31164c5f3306f2d6f1b09a5c6446aa383d4c79/e2e/$mdB"
$zhRow5Url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdB"

$zhXlf6Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/afcfafd7368a4fe5350019f2e3676d53b37ba336/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$zhXlf6"
$zhSrc6Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4f31164c5f3306f2d6f1b09a5c6446aa383d4c79/e2e/$mdC"
$zhRow6Url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdC"

Link-Row $wsZhCn 5 $zhRow5Url $mdB $zhXlf5Url $zhXlf5 $zhSrc5Url $mdB
Link-Row $wsZhCn 6 $zhRow6Url $mdC $zhXlf6Url $zhXlf6 $zhSrc6Url $mdC

# ---------------------------------------------------------------------------
# 4. de-de sheet - rows 5 and 6 (new handoff/handback info)
# ---------------------------------------------------------------------------
$deXlf5 = "2f896ebd-0648-4060-aacb-62692ef7c544.9e4bc652c0f3b7031a37aa768eb56a3aab3d784e.de-de.xlf"
$deXlf6 = "6329960e-3640-4a37-8321-56d9dcf5ee52.2daef6e2b6c383f8d79ffb9a56a7f70dbdb68b3b.de-de.xlf"

Fill-Row $wsDeDe 5 $mdB $deXlf5 "2016-01-25 09:18:49" $mdB "2016-01-25 09:19:45"
Fill-Row $wsDeDe 6 $mdC $deXlf6 "2016-01-25 09:18:49" $mdC "2016-01-25 09:19:45"

$deXlf5Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74283cc9a8ec6c1c113946216d22838362a1e962/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$deXlf5"
$deSrc5Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e95e331e240825b31d89d4f4acc5965a80537a70/e2e/$mdB"
$deRow5Url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdB"

$deXlf6Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74283cc9a8ec6c1c113946216d22838362a1e962/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$deXlf6"
$deSrc6Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e95e331e240825b31d89d4f4acc5965a80537a70/e2e/$mdC"
$deRow6Url = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdC"

Link-Row $wsDeDe 5 $deRow5Url $mdB $deXlf5Url $deXlf5 $deSrc5Url $mdB
Link-Row $wsDeDe 6 $deRow6Url $mdC $deXlf6Url $deXlf6 $deSrc6Url $mdC

Write-Host "Added 3 new rows to Overview, zh-cn and de-de sheets."
